$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.150782
$ws.Range("H2").Value = 3.452345999999999
$ws.Range("I2").Value = 0.03823856951930295
$ws.Range("J2").Value = 0.03823856951930295
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06447966666666667
$ws.Range("N2").Value = 0.193439
$ws.Range("O2").Value = 0.001101138907643723
$ws.Range("P2").Value = 0.001101138907643722
$ws.Range("Q2").Value = 0.07420203976599998
$ws.Range("R2").Value = 0.6678183578939999
$ws.Range("S2").Value = 0.0000421059766703438
$ws.Range("T2").Value = 0.00004210597667034379
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.150782
$ws.Range("H3").Value = 3.452345999999999
$ws.Range("I3").Value = 0.03823856951930295
$ws.Range("J3").Value = 0.03823856951930295
$ws.Range("O3").Value = 0.00657695954769643
$ws.Range("P3").Value = 0.006576959547696431
$ws.Range("Q3").Value = 0.4431991372839999
$ws.Range("R3").Value = 3.988792235555999
$ws.Range("S3").Value = 0.0002514935248902333
$ws.Range("T3").Value = 0.0002514935248902333
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.150782
$ws.Range("H4").Value = 3.452345999999999
$ws.Range("I4").Value = 0.03823856951930295
$ws.Range("J4").Value = 0.03823856951930295
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.008175
$ws.Range("N4").Value = 0.024525
$ws.Range("O4").Value = 0.0001396069650378791
$ws.Range("P4").Value = 0.0001396069650378791
$ws.Range("Q4").Value = 0.009407642849999998
$ws.Range("R4").Value = 0.08466878565
$ws.Range("S4").Value = 0.000005338370637979836
$ws.Range("T4").Value = 0.000005338370637979836
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.150782
$ws.Range("H5").Value = 3.452345999999999
$ws.Range("I5").Value = 0.03823856951930295
$ws.Range("J5").Value = 0.03823856951930295
$ws.Range("M5").Value = 58.099467
$ws.Range("N5").Value = 174.298401
$ws.Range("O5").Value = 0.992182294579622
$ws.Range("P5").Value = 0.992182294579622
$ws.Range("Q5").Value = 66.859820833194
$ws.Range("R5").Value = 601.738387498746
$ws.Range("S5").Value = 0.0379396316471044
$ws.Range("T5").Value = 0.03793963164710439
$ws.Range("I6").Value = 0.9169230158851821
$ws.Range("J6").Value = 0.916923015885182
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.06447966666666667
$ws.Range("N6").Value = 0.193439
$ws.Range("O6").Value = 0.001101138907643723
$ws.Range("P6").Value = 0.001101138907643722
$ws.Range("Q6").Value = 1.779291405049222
$ws.Range("R6").Value = 16.013622645443
$ws.Range("S6").Value = 0.001009659608105197
$ws.Range("T6").Value = 0.001009659608105197
$ws.Range("I7").Value = 0.9169230158851821
$ws.Range("J7").Value = 0.916923015885182
$ws.Range("O7").Value = 0.00657695954769643
$ws.Range("P7").Value = 0.006576959547696431
$ws.Range("S7").Value = 0.006030565583828654
$ws.Range("T7").Value = 0.006030565583828654
$ws.Range("I8").Value = 0.9169230158851821
$ws.Range("J8").Value = 0.916923015885182
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.008175
$ws.Range("N8").Value = 0.024525
$ws.Range("O8").Value = 0.0001396069650378791
$ws.Range("P8").Value = 0.0001396069650378791
$ws.Range("Q8").Value = 0.225585955825
$ws.Range("R8").Value = 2.030273602425
$ws.Range("S8").Value = 0.0001280088394211093
$ws.Range("T8").Value = 0.0001280088394211092
$ws.Range("I9").Value = 0.9169230158851821
$ws.Range("J9").Value = 0.916923015885182
$ws.Range("M9").Value = 58.099467
$ws.Range("N9").Value = 174.298401
$ws.Range("O9").Value = 0.992182294579622
$ws.Range("P9").Value = 0.992182294579622
$ws.Range("Q9").Value = 1603.232268638293
$ws.Range("R9").Value = 14429.09041774464
$ws.Range("S9").Value = 0.9097547818538272
$ws.Range("T9").Value = 0.9097547818538271
$ws.Range("G10").Value = 1.290098666666667
$ws.Range("H10").Value = 3.870296
$ws.Range("I10").Value = 0.04286783035543951
$ws.Range("J10").Value = 0.0428678303554395
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.06447966666666667
$ws.Range("N10").Value = 0.193439
$ws.Range("O10").Value = 0.001101138907643723
$ws.Range("P10").Value = 0.001101138907643722
$ws.Range("Q10").Value = 0.08318513199377778
$ws.Range("R10").Value = 0.748666187944
$ws.Range("S10").Value = 0.00004720343589064507
$ws.Range("T10").Value = 0.00004720343589064505
$ws.Range("G11").Value = 1.290098666666667
$ws.Range("H11").Value = 3.870296
$ws.Range("I11").Value = 0.04286783035543951
$ws.Range("J11").Value = 0.0428678303554395
$ws.Range("O11").Value = 0.00657695954769643
$ws.Range("P11").Value = 0.006576959547696431
$ws.Range("Q11").Value = 0.4968539793617777
$ws.Range("R11").Value = 4.471685814255999
$ws.Range("S11").Value = 0.0002819399861452387
$ws.Range("T11").Value = 0.0002819399861452387
$ws.Range("G12").Value = 1.290098666666667
$ws.Range("H12").Value = 3.870296
$ws.Range("I12").Value = 0.04286783035543951
$ws.Range("J12").Value = 0.0428678303554395
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.008175
$ws.Range("N12").Value = 0.024525
$ws.Range("O12").Value = 0.0001396069650378791
$ws.Range("P12").Value = 0.0001396069650378791
$ws.Range("Q12").Value = 0.0105465566
$ws.Range("R12").Value = 0.0949190094
$ws.Range("S12").Value = 0.000005984647693681575
$ws.Range("T12").Value = 0.000005984647693681574
$ws.Range("G13").Value = 1.290098666666667
$ws.Range("H13").Value = 3.870296
$ws.Range("I13").Value = 0.04286783035543951
$ws.Range("J13").Value = 0.0428678303554395
$ws.Range("M13").Value = 58.099467
$ws.Range("N13").Value = 174.298401
$ws.Range("O13").Value = 0.992182294579622
$ws.Range("P13").Value = 0.992182294579622
$ws.Range("Q13").Value = 74.954044910744
$ws.Range("R13").Value = 674.586404196696
$ws.Range("S13").Value = 0.04253270228570995
$ws.Range("T13").Value = 0.04253270228570994
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.05930433333333333
$ws.Range("H14").Value = 0.177913
$ws.Range("I14").Value = 0.001970584240075516
$ws.Range("J14").Value = 0.001970584240075516
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.06447966666666667
$ws.Range("N14").Value = 0.193439
$ws.Range("O14").Value = 0.001101138907643723
$ws.Range("P14").Value = 0.001101138907643722
$ws.Range("Q14").Value = 0.003823923645222222
$ws.Range("R14").Value = 0.034415312807
$ws.Range("S14").Value = 0.000002169886977536689
$ws.Range("T14").Value = 0.000002169886977536688
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.05930433333333333
$ws.Range("H15").Value = 0.177913
$ws.Range("I15").Value = 0.001970584240075516
$ws.Range("J15").Value = 0.001970584240075516
$ws.Range("O15").Value = 0.00657695954769643
$ws.Range("P15").Value = 0.006576959547696431
$ws.Range("Q15").Value = 0.02283979882422222
$ws.Range("R15").Value = 0.205558189418
$ws.Range("S15").Value = 0.00001296045283230478
$ws.Range("T15").Value = 0.00001296045283230478
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.05930433333333333
$ws.Range("H16").Value = 0.177913
$ws.Range("I16").Value = 0.001970584240075516
$ws.Range("J16").Value = 0.001970584240075516
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.008175
$ws.Range("N16").Value = 0.024525
$ws.Range("O16").Value = 0.0001396069650378791
$ws.Range("P16").Value = 0.0001396069650378791
$ws.Range("Q16").Value = 0.0004848129249999999
$ws.Range("R16").Value = 0.004363316325
$ws.Range("S16").Value = 0.0000002751072851084181
$ws.Range("T16").Value = 0.000000275107285108418
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.05930433333333333
$ws.Range("H17").Value = 0.177913
$ws.Range("I17").Value = 0.001970584240075516
$ws.Range("J17").Value = 0.001970584240075516
$ws.Range("M17").Value = 58.099467
$ws.Range("N17").Value = 174.298401
$ws.Range("O17").Value = 0.992182294579622
$ws.Range("P17").Value = 0.992182294579622
$ws.Range("Q17").Value = 3.445550157457
$ws.Range("R17").Value = 31.009951417113
$ws.Range("S17").Value = 0.001955178792980566
$ws.Range("T17").Value = 0.001955178792980566
